$d = $word.ActiveDocument

# "CIV-6574 removed first paragraph from Witnesses of Fact in Disposal
# Hearing" - the paragraph containing <<disposalHearingWitnessOfFact.input1>>
# ... <<disposalHearingWitnessOfFact.input2>> is removed in its entirety
# (the following paragraphs - input3/input4, input5/input6 - are left
# untouched and simply become the new first/second bullets).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*disposalHearingWitnessOfFact.input1*") {
        $p.Range.Delete()
        break
    }
}
